# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" worksheets, reflecting refreshed scrape counts.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetExhibition.Range("F4").Value = 169
$sheetExhibition.Range("F5").Value = 2794
$sheetExhibition.Range("F7").Value = 221
$sheetExhibition.Range("F8").Value = 20
$sheetExhibition.Range("F10").Value = 68
$sheetExhibition.Range("F11").Value = 78
$sheetExhibition.Range("F12").Value = 2619
$sheetExhibition.Range("F13").Value = 786

$sheetAll = $wb.Worksheets.Item("全部类型")
$sheetAll.Range("F5").Value = 169
$sheetAll.Range("F6").Value = 2794
$sheetAll.Range("F8").Value = 221
$sheetAll.Range("F9").Value = 20
$sheetAll.Range("F12").Value = 68
$sheetAll.Range("F13").Value = 78
$sheetAll.Range("F14").Value = 2619
$sheetAll.Range("F15").Value = 786
